# Add two new analysis columns ("Årsag" and "Ny leverandør") and move the
# previously constant "TCV_range" column out to a new column H, while
# repurposing column F to hold the cancellation reason ("Årsag").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# Match the header formatting (bold font, border, centered) already used
# by the other header cells.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Move the constant TCV_range value from F to the new H column -------
$ws.Range("H2:H12").Value = "80000-100000"

# --- Fill in the cancellation reason in column F -------------------------
$ws.Range("F2").Value = "Ikke oplyst"
$ws.Range("F3").Value = "Ikke oplyst"
$ws.Range("F4").Value = "Ikke oplyst"
$ws.Range("F5").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F6").Value = "Ikke oplyst"
$ws.Range("F7").Value = "Virksomheden lukker"
$ws.Range("F8").Value = "Ikke flere medarbejdere i virksomheden"
$ws.Range("F9").Value = "Virksomheden lukker"
$ws.Range("F10").Value = "Ikke oplyst"
$ws.Range("F11").Value = "Virksomheden lukker"
$ws.Range("F12").Value = "Insourcing af lønnen (anden leverandør)"

# --- Fill in the new supplier in column G where known --------------------
$ws.Range("G2").Value = "Intect"
$ws.Range("G6").Value = "Intect"
$ws.Range("G10").Value = "Zenegy"

$ws.Range("A1").Select()
